$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '68.306.12'
$ws.Range('E2').Value = '  +3.23%  '

# Row 3
Set-TextValue 'D3' '3.643.91'
$ws.Range('E3').Value = '  +2.82%  '

# Row 4
$ws.Range('E4').Value = '  -0.20%  '

# Row 5
Set-TextValue 'D5' '202.12'
$ws.Range('E5').Value = '  +12.40%  '

# Row 6
Set-TextValue 'D6' '574.42'
$ws.Range('E6').Value = '  -1.82%  '

# Row 7
Set-TextValue 'D7' '3.637.12'
$ws.Range('E7').Value = '  +2.97%  '

# Row 8
Set-TextValue 'D8' '0.619'

# Row 9
$ws.Range('E9').Value = '  -0.43%  '

# Row 10
Set-TextValue 'D10' '0.679'
$ws.Range('E10').Value = '  +2.70%  '

# Row 11
Set-TextValue 'D11' '0.154'
$ws.Range('E11').Value = '  +9.15%  '

# Row 12
Set-TextValue 'D12' '57.50'
$ws.Range('E12').Value = '  +8.66%  '

# Row 13
Set-TextValue 'D13' '0.0000294'
$ws.Range('E13').Value = '  +19.74%  '

# Row 14
Set-TextValue 'D14' '10.09'
$ws.Range('E14').Value = '  +4.54%  '

# Row 15
Set-TextValue 'D15' '4.228.30'
$ws.Range('E15').Value = '  +2.25%  '

# Row 16
Set-TextValue 'D16' '3.651.74'
$ws.Range('E16').Value = '  +2.80%  '

# Row 17
$ws.Range('E17').Value = '  +0.81%  '

# Row 18
Set-TextValue 'D18' '12.49'
$ws.Range('E18').Value = '  +4.22%  '

# Row 19
Set-TextValue 'D19' '68.272.48'
$ws.Range('E19').Value = '  +3.42%  '

# Row 20
Set-TextValue 'D20' '18.65'
$ws.Range('E20').Value = '  +2.72%  '

# Row 21
$ws.Range('E21').Value = '  +4.57%  '

# Row 22
Set-TextValue 'D22' '403.25'
$ws.Range('E22').Value = '  +3.73%  '

# Row 23
Set-TextValue 'D23' '13.05'
$ws.Range('E23').Value = '  +29.25%  '

# Row 24
Set-TextValue 'D24' '4.22'
$ws.Range('E24').Value = '  -0.87%  '

# Row 25
Set-TextValue 'D25' '86.13'
$ws.Range('E25').Value = '  +2.12%  '

# Row 26
Set-TextValue 'D26' '2.97'
$ws.Range('E26').Value = '  +5.11%  '

# Row 27
Set-TextValue 'D27' '12.63'
$ws.Range('E27').Value = '  +3.93%  '

# Row 28
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D28' '6.11'
$ws.Range('E28').Value = '  +1.69%  '

# Row 29
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D29' '3.84'
$ws.Range('E29').Value = '  +8.99%  '

# Row 30
Set-TextValue 'D30' '8.27'
$ws.Range('E30').Value = '  +24.66%  '

# Row 31
Set-TextValue 'D31' '9.18'
$ws.Range('E31').Value = '  +4.48%  '

# Row 32
Set-TextValue 'D32' '32.02'
$ws.Range('E32').Value = '  +4.49%  '

# Row 33
Set-TextValue 'D33' '697.51'
$ws.Range('E33').Value = '  +14.58%  '

# Row 34
Set-TextValue 'D34' '12.26'
$ws.Range('E34').Value = '  +3.50%  '

# Row 35
$ws.Range('E35').Value = '  +5.59%  '

# Row 36
Set-TextValue 'D36' '64.37'
$ws.Range('E36').Value = '  -0.98%  '

# Row 37
Set-TextValue 'D37' '42.79'
$ws.Range('E37').Value = '  +4.75%  '

# Row 38
Set-TextValue 'D38' '0.429'
$ws.Range('E38').Value = '  +16.92%  '

# Row 39
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D39' '1.00'
$ws.Range('E39').Value = '  +0.08%  '

# Row 40
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D40' '0.0₃0796'
$ws.Range('E40').Value = '  +8.95%  '

# Row 41
Set-TextValue 'D41' '0.141'
$ws.Range('E41').Value = '  +10.06%  '

# Row 42
Set-TextValue 'D42' '3.252.84'
$ws.Range('E42').Value = '  +13.95%  '

# Row 43
Set-TextValue 'D43' '3.15'
$ws.Range('E43').Value = '  +14.42%  '

# Row 44
Set-TextValue 'D44' '2.82'
$ws.Range('E44').Value = '  +18.83%  '

# Row 45
Set-TextValue 'D45' '0.998'
$ws.Range('E45').Value = '  -0.27%  '

# Row 46
$ws.Range('E46').Value = '  +39.87%  '

# Row 47
Set-TextValue 'D47' '0.0421'
$ws.Range('E47').Value = '  +4.54%  '

# Row 48
Set-TextValue 'D48' '8.94'
$ws.Range('E48').Value = '  +9.36%  '

# Row 49
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D49' '2.69'
$ws.Range('E49').Value = '  +9.40%  '

# Row 50
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D50' '0.132'
$ws.Range('E50').Value = '  +2.56%  '

# Row 51
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D51' '142.00'
$ws.Range('E51').Value = '  +3.71%  '
